$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.886.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.815.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.91%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9943'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.80%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9938'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3938'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3491'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.28'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.209'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07569'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9920'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.544'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.810.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.197'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001108'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06689'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9940'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.592'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.847.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.95%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.398'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.564'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.497'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.95%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '155.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.013.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '135.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.027'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.68%  '

$ws.Range("E33").Value = '  -0.69%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08827'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.33'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.546'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02427'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6927'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06539'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.613'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2228'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.265'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.568'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.66'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.76%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6554'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.79%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9937'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.862'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.167'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.43'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07244'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.37%  '
